# update matrix mult seq with one add one mult
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10: move the PCF-tool numbers that used to live in L:Q into the
# Synopsys D:J columns (with a multiplier column I10), then blank out L:Q
# (keeping their number formatting).
# ---------------------------------------------------------------------------
$ws.Range("D10").Value2 = 2047
$ws.Range("E10").Value2 = 2050
$ws.Range("D10:E10").NumberFormat = "#,##0"
$ws.Range("F10").Formula = "=SUM(D10:E10)"
$ws.Range("F10").NumberFormat = "#,##0"
$ws.Range("G10").Value2 = 3070
$ws.Range("G10").NumberFormat = "#,##0"
$ws.Range("H10").Formula = "=SUM(D10:F10)"
$ws.Range("H10").NumberFormat = "#,##0"
$ws.Range("I10").Value2 = 1
$ws.Range("J10").Formula = "=(SUM(D10:E10)+5*G10)*I10"
$ws.Range("J10").NumberFormat = "#,##0"
$ws.Range("J10").Font.Bold = $true

$ws.Range("L10:Q10").ClearContents()

# ---------------------------------------------------------------------------
# Rows 19, 20, 21, 22: add explicit (empty) style-0 cells in columns A/B,
# overriding the bold column default formatting used elsewhere in those
# columns.
# ---------------------------------------------------------------------------
foreach ($r in 19, 20) {
    $cell = $ws.Range("A$r")
    $cell.NumberFormat = "GENERAL"
    $cell.Font.Bold = $false
}
foreach ($r in 21, 22) {
    $cellA = $ws.Range("A$r")
    $cellA.NumberFormat = "GENERAL"
    $cellA.Font.Bold = $false
    $cellB = $ws.Range("B$r")
    $cellB.NumberFormat = "GENERAL"
    $cellB.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# Row 23: was the "3x3" MatmultSeq row, now becomes the "2x2" row with new
# measured values and a computed (2*2*2) multiplier.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value2 = "2x2"
$ws.Range("E23").Value2 = 997
$ws.Range("G23").Value2 = 1956
$ws.Range("I23").Formula = "=2*2*2"

# ---------------------------------------------------------------------------
# New row 24: "3x3" MatmultSeq row.
# ---------------------------------------------------------------------------
$ws.Range("B24").Value2 = "3x3"
$ws.Range("B24").Font.Bold = $true
$ws.Range("B24").NumberFormat = "GENERAL"

$ws.Range("D24").Value2 = 1026
$ws.Range("E24").Value2 = 997
$ws.Range("F24").Formula = "=SUM(D24:E24)"
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("G24").Value2 = 1956
$ws.Range("H24").Formula = "=SUM(D24:F24)"
$ws.Range("H24").NumberFormat = "#,##0"
$ws.Range("I24").Formula = "=3*3*3"
$ws.Range("J24").Formula = "=(SUM(D24:E24)+5*G24)*I24"
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("J24").Font.Bold = $true

# ---------------------------------------------------------------------------
# New row 25: "5x5" MatmultSeq row.
# ---------------------------------------------------------------------------
$ws.Range("B25").Value2 = "5x5"
$ws.Range("B25").Font.Bold = $true
$ws.Range("B25").NumberFormat = "GENERAL"

$ws.Range("D25").Value2 = 1026
$ws.Range("E25").Value2 = 997
$ws.Range("F25").Formula = "=SUM(D25:E25)"
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("G25").Value2 = 1956
$ws.Range("H25").Formula = "=SUM(D25:F25)"
$ws.Range("H25").NumberFormat = "#,##0"
$ws.Range("I25").Formula = "=5*5*5"
$ws.Range("J25").Formula = "=(SUM(D25:E25)+5*G25)*I25"
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("J25").Font.Bold = $true

# ---------------------------------------------------------------------------
# Update the view: scroll so column G is the left-most visible column, and
# select M12.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("M12").Select()
